# Insert a new data row before existing row 119 (shifting all following
# rows down by one) and populate it with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(119).Insert()

$ws.Cells.Item(119, 1).Value2  = 4
$ws.Cells.Item(119, 2).Value2  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(119, 3).Value2  = "Los Lagos"
$ws.Cells.Item(119, 4).Value2  = 44510
$ws.Cells.Item(119, 5).Value2  = 10
$ws.Cells.Item(119, 6).Value2  = 100112043
$ws.Cells.Item(119, 7).Value2  = "Pepino ensalada"
$ws.Cells.Item(119, 8).Value2  = "Sin especificar"
$ws.Cells.Item(119, 9).Value2  = "Primera"
$ws.Cells.Item(119, 10).Value2 = 100
$ws.Cells.Item(119, 11).Value2 = 12000
$ws.Cells.Item(119, 12).Value2 = 12000
$ws.Cells.Item(119, 13).Value2 = 12000
$ws.Cells.Item(119, 14).Value2 = "`$/caja 60 unidades"
$ws.Cells.Item(119, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(119, 16).Value2 = 200
$ws.Cells.Item(119, 17).Value2 = 60
$ws.Cells.Item(119, 18).Value2 = "Hortaliza"
